$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row before row 41, shifting rows 41:48 down to 42:49
$ws.Rows.Item(41).Insert()

# Fill the new row 41
$ws.Range("A41").Value = "juenger als 23 oder vor 1940 geboren?"

# Match number formatting of B40 (integer) on the new B41 cell
$ws.Range("B40").Copy()
$ws.Range("B41").PasteSpecial(-4122)  # xlPasteFormats

# Update selection / view to match target state
$ws.Range("B40").Select()
